$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1) Remove the "Meta description: ..." paragraph that currently sits
#    right after the H1 title paragraph.
# ------------------------------------------------------------------
$metaPara = $d.Paragraphs.Item(2)
$metaPara.Range.Delete()

# ------------------------------------------------------------------
# 2) Insert a new bold paragraph ("Play Double Triple Chance for Free
#    - Exciting Fruit Machine Slot") right before the final paragraph
#    (the one that used to hold the "Create a feature image..." prompt).
#    We build it by duplicating a plain (non-bold / non-italic) body
#    paragraph so the new paragraph break is created cleanly, then we
#    overwrite its text and apply bold formatting.
# ------------------------------------------------------------------
$srcPara = $d.Paragraphs.Item(4)
$srcParaEnd = $d.Paragraphs.Item(5).Range.Start
$fullSrcRange = $d.Range($srcPara.Range.Start, $srcParaEnd)
$fullSrcRange.Copy()

$count = $d.Paragraphs.Count
$last = $d.Paragraphs.Item($count)
$insertPoint = $d.Range($last.Range.Start, $last.Range.Start)
$insertPoint.Paste()

$newPara = $d.Paragraphs.Item($count)
$newRange = $d.Range($newPara.Range.Start, $newPara.Range.End)
$newRange.Text = "Play Double Triple Chance for Free - Exciting Fruit Machine Slot"
$newRange.Font.Bold = 1

# ------------------------------------------------------------------
# 3) Replace the text of the (now shifted) final paragraph - the old
#    "Create a feature image..." image prompt - with the meta
#    description copy, keeping its existing italic formatting intact.
# ------------------------------------------------------------------
$finalCount = $d.Paragraphs.Count
$finalPara = $d.Paragraphs.Item($finalCount)
$finalRange = $d.Range($finalPara.Range.Start, $finalPara.Range.End)
$finalRange.Text = "Get ready for a dynamic online slot experience with Double Triple Chance. Play for free and win big with Bonus Rewin feature and wide betting options."
